# Apply updated T20 ("E" column) appearance counts for several players,
# reflecting newly-scraped batting performance data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_formats_raw")

$updates = @{
    14 = 41
    15 = 8
    17 = 45
    31 = 51
    34 = 12
    35 = 35
    39 = 24
    50 = 75
    52 = 58
    54 = 3
    58 = 26
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
